# Add five new attraction rows to the locations sheet. The existing three
# data rows are re-typed too, which is what drops the sheet's assorted
# one-off cell styles (every data cell but the header-style A2 ends up
# back on the plain default style), and the selection is left at F4 to
# match the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing data rows (2-4) entirely - this drops their cached
# per-cell styles *and* their stored row heights, rather than merely
# blanking the cells, so the rows we retype below come back on defaults.
$ws.Rows("2:4").Delete()

# Fill in the "attraction" column first, in the same order the names were
# originally typed, so the shared-string table is rebuilt with that exact
# ordering (the three existing names simply reuse their old string
# entries; the five new ones are appended in entry order, not row order).
$ws.Cells.Item(2, 1).Value = "Despicable Me Minion Mayhem"
$ws.Cells.Item(3, 1).Value = "Hollywood Rip Ride Rockit"
$ws.Cells.Item(4, 1).Value = "Universal Music Plaza Stage"
$ws.Cells.Item(6, 1).Value = "Race Through New York Starring Jimmy Fallon"
$ws.Cells.Item(7, 1).Value = "Harry Potter and the Escape from Gringotts"
$ws.Cells.Item(5, 1).Value = "Hogwarts Express"
$ws.Cells.Item(8, 1).Value = "Men in Black: Alien Attack"
$ws.Cells.Item(9, 1).Value = "The Simpsons Ride"

# Now fill in the lat/lon columns for every row.
$ws.Cells.Item(2, 2).Value = 28.475272
$ws.Cells.Item(2, 3).Value = -81.468102999999999

$ws.Cells.Item(3, 2).Value = 28.474900000000002
$ws.Cells.Item(3, 3).Value = -81.468299999999999

$ws.Cells.Item(4, 2).Value = 28.475477999999999
$ws.Cells.Item(4, 3).Value = -81.468857999999997

$ws.Cells.Item(5, 2).Value = 28.479399999999998
$ws.Cells.Item(5, 3).Value = -81.470299999999995

$ws.Cells.Item(6, 2).Value = 28.475683
$ws.Cells.Item(6, 3).Value = -81.469449999999995

$ws.Cells.Item(7, 2).Value = 28.480277999999998
$ws.Cells.Item(7, 3).Value = -81.47

$ws.Cells.Item(8, 2).Value = 28.480858000000001
$ws.Cells.Item(8, 3).Value = -81.467399999999998

$ws.Cells.Item(9, 2).Value = 28.479438999999999
$ws.Cells.Item(9, 3).Value = -81.467364000000003

# Drop the leftover per-column width/style formatting on B:C (column A
# keeps its bestFit width) carried over from the old sheet.
$ws.Columns("A:C").ClearFormats()

# Re-apply the one remaining bit of manual formatting: the first
# attraction's name is shown a little larger and in a dark grey.
$ws.Range("A2").Font.Size = 14
$ws.Range("A2").Font.Color = 2236704

# Match the saved selection / scroll position.
$ws.Range("F4").Select()
